$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5.981099999999998
$ws.Range("A12").Value = -21.40019999999999
$ws.Range("B23").Value = 8.950099999999996
$ws.Range("B28").Value = 5.447200000000002
$ws.Range("A32").Value = -21.04559999999998
$ws.Range("B32").Value = 5.965399999999995
$ws.Range("B34").Value = 9.811100000000007
$ws.Range("A36").Value = -20.1445
$ws.Range("A38").Value = -19.94759999999998
$ws.Range("B42").Value = 9.964299999999996
$ws.Range("A46").Value = -21.91320000000001
$ws.Range("A54").Value = -21.9908
$ws.Range("B54").Value = 4.912899999999997
$ws.Range("A55").Value = -22.0502
$ws.Range("A67").Value = -21.40759999999997
$ws.Range("A69").Value = -21.52769999999998
$ws.Range("A72").Value = -22.0486
$ws.Range("A91").Value = -20.66849999999998
$ws.Range("B97").Value = 5.762899999999997
$ws.Range("A99").Value = -21.95009999999999
$ws.Range("B99").Value = 5.540499999999997
$ws.Range("B101").Value = 4.450699999999999
$ws.Range("A104").Value = -21.5665
